# Updated cryptos list on Tue Feb  6 02:51:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.884.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.328.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.55%  '

$ws.Range("E7").Value = '  +0.89%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.72'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.24'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.35%  '

$ws.Range("E12").Value = '  +1.08%  '

$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.688.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.332.43'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.28%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.808.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0892'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.10'
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.30'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.22%  '

$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("E28").Value = '  -0.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.52'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.59%  '

$ws.Range("E32").Value = '  -0.07%  '

$ws.Range("E33").Value = '  +2.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.99%  '

$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("E36").Value = '  +2.95%  '

$ws.Range("E37").Value = '  -1.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.79'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.99%  '

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.84%  '

$ws.Range("E41").Value = '  +0.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.953.59'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  +1.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.32%  '

$ws.Range("E46").Value = '  +3.41%  '

$ws.Range("E47").Value = '  +0.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.555.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.52'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.75%  '
